$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the duplicate "poverty_level_total" row (old row 129).
#    Everything below shifts up by one (old 130 -> new 129, etc).
# ------------------------------------------------------------------
$ws.Rows(129).Delete()

# ------------------------------------------------------------------
# 2. Append the new B08124 "occupation" rows at the bottom of the
#    table (new rows 157-163). Values are written in an order that
#    reproduces the original authoring sequence for the shared
#    string table (some label columns were filled out of strict
#    row order, e.g. col B for rows 159/158, then col D for row 157).
# ------------------------------------------------------------------
$ws.Range("A157").Value = "B08124_001E"
$ws.Range("B157").Value = "total_occupation"
$ws.Range("A158").Value = "B08124_002E"
$ws.Range("A159").Value = "B08124_003E"
$ws.Range("B159").Value = "service_occupation"
$ws.Range("B158").Value = "management_business_occupation"
$ws.Range("D157").Value = "occupation"
$ws.Range("A160").Value = "B08124_004E"
$ws.Range("B160").Value = "sales_office_occupation"
$ws.Range("A161").Value = "B08124_005E"
$ws.Range("B161").Value = "natural_res_construction_occupation"
$ws.Range("A162").Value = "B08124_006E"
$ws.Range("B162").Value = "production_transp_occupation"
$ws.Range("A163").Value = "B08124_007E"
$ws.Range("B163").Value = "military_occupation"

$ws.Range("C157").Value = "total_occupation"
$ws.Range("C158").Value = "total_occupation"
$ws.Range("D158").Value = "occupation"
$ws.Range("C159").Value = "total_occupation"
$ws.Range("D159").Value = "occupation"
$ws.Range("C160").Value = "total_occupation"
$ws.Range("D160").Value = "occupation"
$ws.Range("C161").Value = "total_occupation"
$ws.Range("D161").Value = "occupation"
$ws.Range("C162").Value = "total_occupation"
$ws.Range("D162").Value = "occupation"
$ws.Range("C163").Value = "total_occupation"
$ws.Range("D163").Value = "occupation"

# ------------------------------------------------------------------
# 3. Update the window scroll position / selection to match where
#    the author left the cursor after the edit.
# ------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 107
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A129:XFD129").Select()
